$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "34.490.68"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.804.11"
$ws.Range("E3").Value = "  +0.25%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.11%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "

# Row 6: XRP
$ws.Range("E6").Value = "  +3.98%  "

# Row 8: Solana
$ws.Range("E8").Value = "  +6.15%  "

# Row 9: Cardano
$ws.Range("E9").Value = "  +1.23%  "

# Row 10: Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0694"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.52%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0954"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.32%  "

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.064.61"
$ws.Range("E12").Value = "  +0.28%  "

# Row 13: Chainlink
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "

# Row 14: WrappedEther
$ws.Range("D14").Value = "1.801.34"
$ws.Range("E14").Value = "  +0.05%  "

# Row 15: Polygon
$ws.Range("E15").Value = "  +0.60%  "

# Row 16: WrappedBTC
$ws.Range("D16").Value = "34.460.82"

# Row 17: Polkadot
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.37"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.49%  "

# Row 18: Litecoin
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.11%  "

# Row 19: ShibaInu
$ws.Range("D19").Value = "0.0₃0798"
$ws.Range("E19").Value = "  -0.80%  "

# Row 20: BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.84%  "

# Row 21: Avalanche
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.69%  "

# Row 22: Dai
$ws.Range("E22").Value = "  +0.06%  "

# Row 23: Uniswap
$ws.Range("E23").Value = "  -0.24%  "

# Row 24: Monero
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "173.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.52%  "

# Row 25: Toncoin
$ws.Range("E25").Value = "  +2.12%  "

# Row 26: Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.85"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.48%  "

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.30%  "

# Row 28: Stellar
$ws.Range("E28").Value = "  +2.06%  "

# Row 29: BinanceUSD
$ws.Range("E29").Value = "  -0.08%  "

# Row 30: InternetComputer(DFINITY)
$ws.Range("E30").Value = "  -1.92%  "

# Row 31: Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0532"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.28%  "

# Row 32: 'Filecoin' -> 'PancakeSwap'
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.56%  "

# Row 33: 'PancakeSwap' -> 'Filecoin'
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.75%  "

# Row 34: LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.32%  "

# Row 35: ImmutableX
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.686"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.03%  "

# Row 36: Maker
$ws.Range("D36").Value = "1.394.58"
$ws.Range("E36").Value = "  -2.66%  "

# Row 37: RenderToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.37%  "

# Row 38: TrustWalletToken
$ws.Range("E38").Value = "  -1.04%  "

# Row 39: VeChain
$ws.Range("E39").Value = "  -1.15%  "

# Row 40: Aave
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.43%  "

# Row 41: MXToken
$ws.Range("E41").Value = "  +2.40%  "

# Row 42: ARBITRUM
$ws.Range("E42").Value = "  +1.50%  "

# Row 43: HuobiToken
$ws.Range("E43").Value = "  -0.84%  "

# Row 44: InjectiveProtocol
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.95%  "

# Row 45: WEMIXToken
$ws.Range("E45").Value = "  +3.58%  "

# Row 46: Kaspa
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0511"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.14%  "

# Row 47: FraxShare
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.10%  "

# Row 48: RocketPoolETH
$ws.Range("D48").Value = "1.964.62"
$ws.Range("E48").Value = "  +0.34%  "

# Row 49: Quant
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.17%  "

# Row 50: PaxDollar
$ws.Range("E50").Value = "  +0.06%  "

# Row 51: BabyDogeCoin
$ws.Range("E51").Value = "  +1.52%  "
